$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Secondary Users: Support users" -> "Secondary Users: Support team
# members" (paragraph "Secondary Users: Support users").
# We only want the final, standalone "users" run's text to change; the
# preceding " Support " run must stay untouched / stay a separate run.
# A Bookmark is inserted at the run boundary purely to force the text-edit
# engine to keep the runs split (it gets removed again immediately after);
# it leaves no residue in the saved document.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Secondary Users:") -and $t.Contains("Support users")) {
        $pStart = $p.Range.Start
        $idx = $t.LastIndexOf("users")
        $wordStart = $pStart + $idx
        $wordEnd = $wordStart + 5

        $d.Bookmarks.Add("zz_split_a", $d.Range($wordStart, $wordStart))
        $d.Range($wordStart, $wordEnd).Text = "team members"
        $d.Bookmarks("zz_split_a").Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Edit 2: "...assign different tags to different support users..." ->
# "...assign different tags to different support team members...", keeping
# " so that" as its own (unchanged) italic run. The original text is a
# single run "<nbsp>assign different tags to different support users<nbsp>".
# The target splits it into three runs:
#   " assign different tags to different support "  (regular spaces now)
#   "team members"
#   " "                                              (regular space now)
# Two temporary bookmarks force the run boundaries; they are removed right
# after, leaving three plain runs with no leftover formatting.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Contains("assign different tags to different support")) {
        $pStart = $p.Range.Start
        $idx = $t.IndexOf("assign different tags to different support users")
        $runStart = $pStart + $idx - 1
        $runEnd = $pStart + $idx + ("assign different tags to different support users").Length + 1
        $wordStart = $runEnd - 6
        $wordEnd = $runEnd - 1

        $d.Bookmarks.Add("zz_split_b1", $d.Range($wordStart, $wordStart))
        $d.Bookmarks.Add("zz_split_b2", $d.Range($wordEnd, $wordEnd))

        $d.Range($runStart, $wordStart).Text = " assign different tags to different support "
        $d.Range($wordStart, $wordEnd).Text = "team members"
        $newWordEnd = $wordEnd + 7
        $d.Range($newWordEnd, $newWordEnd + 1).Text = " "

        $d.Bookmarks("zz_split_b1").Delete()
        $d.Bookmarks("zz_split_b2").Delete()
        break
    }
}
